$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 231. This shifts the existing rows
# 231-290 down to 232-291 (Excel's InsertRow semantics), matching the
# diff where row 290 becomes 291 and so on, with a brand-new record
# landing on row 231.
$ws.Rows.Item(231).Insert()

# Populate the newly inserted row 231 with the new record's data.
$ws.Range("A231").Value = 5
$ws.Range("B231").Value = "Macroferia Regional de Talca"
$ws.Range("C231").Value = "Maule"
$ws.Range("D231").Value = 44782
$ws.Range("E231").Value = 7
$ws.Range("F231").Value = 100112009
$ws.Range("G231").Value = "Acelga"
$ws.Range("H231").Value = "Sin especificar"
$ws.Range("I231").Value = "Primera"
$ws.Range("J231").Value = 500
$ws.Range("K231").Value = 3000
$ws.Range("L231").Value = 3000
$ws.Range("M231").Value = 3000
$ws.Range("N231").Value = "$/docena de atados (4 kilos)"
$ws.Range("O231").Value = "Región del Maule"
$ws.Range("P231").Value = 750
$ws.Range("Q231").Value = 4
$ws.Range("R231").Value = "Hortaliza"
